$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.610.07"
$ws.Range("E2").Value = "  +0.54%  "

$ws.Range("D3").Value = "1.636.21"
$ws.Range("E3").Value = "  -0.53%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'212.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.05%  "

$ws.Range("E6").Value = "  -1.01%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").Value = "'22.95"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.76%  "

$ws.Range("E9").Value = "  +0.13%  "

$ws.Range("E10").Value = "  -0.17%  "

$ws.Range("D11").Value = "'0.0893"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.17%  "

$ws.Range("D12").Value = "1.868.28"
$ws.Range("E12").Value = "  -0.47%  "

$ws.Range("D13").Value = "1.631.11"
$ws.Range("E13").Value = "  +0.28%  "

$ws.Range("E14").Value = "  -0.16%  "

$ws.Range("D15").Value = "'0.558"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.58%  "

$ws.Range("D16").Value = "'64.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.13%  "

$ws.Range("D17").Value = "27.606.78"
$ws.Range("E17").Value = "  +0.68%  "

$ws.Range("D18").Value = "'228.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.70%  "

$ws.Range("D19").Value = "'7.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.71%  "

$ws.Range("D20").Value = "0.0₃0721"
$ws.Range("E20").Value = "  -0.21%  "

$ws.Range("E21").Value = "  +0.10%  "

$ws.Range("E22").Value = "  -1.08%  "

$ws.Range("D23").Value = "'10.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.89%  "

$ws.Range("D24").Value = "'1.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.66%  "

$ws.Range("D25").Value = "'150.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.05%  "

$ws.Range("E26").Value = "  -1.33%  "

$ws.Range("E27").Value = "  -1.86%  "

$ws.Range("E28").Value = "  +0.01%  "

$ws.Range("D29").Value = "'15.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.29%  "

$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("D31").Value = "'0.0485"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.04%  "

$ws.Range("D32").Value = "'3.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.19%  "

$ws.Range("D33").Value = "1.454.93"
$ws.Range("E33").Value = "  +2.28%  "

$ws.Range("D34").Value = "'3.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.97%  "

$ws.Range("E35").Value = "  -1.03%  "

$ws.Range("E36").Value = "  -0.27%  "

$ws.Range("D37").Value = "'0.564"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.45%  "

$ws.Range("D38").Value = "'0.876"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.97%  "

$ws.Range("E39").Value = "  +0.38%  "

$ws.Range("D40").Value = "'0.893"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.68%  "

$ws.Range("D41").Value = "'69.67"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.70%  "

$ws.Range("E42").Value = "  +0.08%  "

$ws.Range("E43").Value = "  -1.04%  "

$ws.Range("E44").Value = "  +1.56%  "

$ws.Range("B45").Value = "MXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D45").Value = "'2.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.73%  "

$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.778.23"
$ws.Range("E46").Value = "  -0.47%  "

$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "'1.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.49%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'86.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.08%  "

$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0107"
$ws.Range("E49").Value = "  -0.56%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.0984"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.03%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.48%  "
